# Scheduled runner update: refresh market-price-driven leve profit figures
# across the per-job sheets (currentAveragePrice* / LevePrice* / LeveProfit* columns).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H30").Value = 5050036
$ws.Range("I30").Value = 5050036
$ws.Range("K30").Value = 15150108
$ws.Range("M30").Value = -15150007
$ws.Range("H62").Value = 3803.5789
$ws.Range("I62").Value = 3579.3125
$ws.Range("K62").Value = 3579.3125
$ws.Range("M62").Value = -2955.3125
$ws.Range("H65").Value = 3803.5789
$ws.Range("I65").Value = 3579.3125
$ws.Range("K65").Value = 17896.5625
$ws.Range("M65").Value = -14776.5625
$ws.Range("H98").Value = 2871.423
$ws.Range("I98").Value = 1973.2084
$ws.Range("J98").Value = 13650
$ws.Range("K98").Value = 1973.2084
$ws.Range("L98").Value = 13650
$ws.Range("M98").Value = -475.2084
$ws.Range("N98").Value = -16646
$ws.Range("H113").Value = 23497.46
$ws.Range("I113").Value = 54984
$ws.Range("K113").Value = 54984
$ws.Range("M113").Value = -51730
$ws.Range("H116").Value = 5443.375
$ws.Range("J116").Value = 5634
$ws.Range("L116").Value = 5634
$ws.Range("N116").Value = -12518
$ws.Range("H122").Value = 2871.423
$ws.Range("I122").Value = 1973.2084
$ws.Range("J122").Value = 13650
$ws.Range("K122").Value = 5919.6252
$ws.Range("L122").Value = 40950
$ws.Range("M122").Value = -3469.6252
$ws.Range("N122").Value = -45850
$ws.Range("H137").Value = 946.94446
$ws.Range("I137").Value = 879.36365
$ws.Range("J137").Value = 1053.1428
$ws.Range("K137").Value = 2638.09095
$ws.Range("L137").Value = 3159.4284
$ws.Range("M137").Value = -88.09094999999979
$ws.Range("N137").Value = -8259.428400000001
$ws.Range("H138").Value = 5548.5654
$ws.Range("J138").Value = 7420.6206
$ws.Range("L138").Value = 22261.8618
$ws.Range("N138").Value = -32541.8618

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3325.7778
$ws.Range("I32").Value = 1991.4348
$ws.Range("K32").Value = 1991.4348
$ws.Range("M32").Value = -1704.4348
$ws.Range("H45").Value = 1617.0834
$ws.Range("J45").Value = 1632.6666
$ws.Range("L45").Value = 1632.6666
$ws.Range("N45").Value = -2386.6666
$ws.Range("H74").Value = 1856.3448
$ws.Range("I74").Value = 1289.6471
$ws.Range("K74").Value = 1289.6471
$ws.Range("M74").Value = -415.6470999999999
$ws.Range("H77").Value = 1856.3448
$ws.Range("I77").Value = 1289.6471
$ws.Range("K77").Value = 6448.2355
$ws.Range("M77").Value = -2080.2355
$ws.Range("H122").Value = 6160.409
$ws.Range("I122").Value = 6011.5757
$ws.Range("K122").Value = 18034.7271
$ws.Range("M122").Value = -15584.7271
$ws.Range("H132").Value = 7443.9395
$ws.Range("I132").Value = 7715.9463
$ws.Range("K132").Value = 23147.8389
$ws.Range("M132").Value = -20617.8389

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H54").Value = 21353
$ws.Range("I54").Value = 20530
$ws.Range("J54").Value = 22999
$ws.Range("K54").Value = 20530
$ws.Range("L54").Value = 22999
$ws.Range("M54").Value = -20046
$ws.Range("N54").Value = -23967
$ws.Range("H99").Value = 1966.8948
$ws.Range("I99").Value = 2015.0555
$ws.Range("K99").Value = 2015.0555
$ws.Range("M99").Value = -517.0554999999999
$ws.Range("H134").Value = 2017.6154
$ws.Range("I134").Value = 2006.238
$ws.Range("K134").Value = 6018.714
$ws.Range("M134").Value = -3483.714

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 3979.5
$ws.Range("I62").Value = 3759.4285
$ws.Range("K62").Value = 3759.4285
$ws.Range("M62").Value = -3135.4285
$ws.Range("H65").Value = 3979.5
$ws.Range("I65").Value = 3759.4285
$ws.Range("K65").Value = 18797.1425
$ws.Range("M65").Value = -15677.1425
$ws.Range("H99").Value = 1406.5555
$ws.Range("I99").Value = 1406.5555
$ws.Range("K99").Value = 1406.5555
$ws.Range("M99").Value = 91.44450000000006
$ws.Range("H126").Value = 1406.5555
$ws.Range("I126").Value = 1406.5555
$ws.Range("K126").Value = 4219.666499999999
$ws.Range("M126").Value = -1749.666499999999
$ws.Range("H134").Value = 1960
$ws.Range("I134").Value = 1692.091
$ws.Range("J134").Value = 3223
$ws.Range("K134").Value = 5076.272999999999
$ws.Range("L134").Value = 9669
$ws.Range("M134").Value = -2541.272999999999
$ws.Range("N134").Value = -14739

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 2435.5715
$ws.Range("I3").Value = 2435.5715
$ws.Range("K3").Value = 7306.7145
$ws.Range("M3").Value = -7194.7145
$ws.Range("H4").Value = 28682074
$ws.Range("I4").Value = 35848104
$ws.Range("K4").Value = 107544312
$ws.Range("M4").Value = -107544200
$ws.Range("H5").Value = 866.61536
$ws.Range("J5").Value = 1153.125
$ws.Range("L5").Value = 3459.375
$ws.Range("N5").Value = -3683.375
$ws.Range("H40").Value = 255.26315
$ws.Range("I40").Value = 115.625
$ws.Range("J40").Value = 1000
$ws.Range("K40").Value = 462.5
$ws.Range("L40").Value = 4000
$ws.Range("M40").Value = -393.5
$ws.Range("N40").Value = -4138
$ws.Range("H68").Value = 1909.1875
$ws.Range("I68").Value = 1239.2858
$ws.Range("J68").Value = 2430.2222
$ws.Range("K68").Value = 3717.8574
$ws.Range("L68").Value = 7290.6666
$ws.Range("M68").Value = -2906.8574
$ws.Range("N68").Value = -8912.6666
$ws.Range("H71").Value = 1909.1875
$ws.Range("I71").Value = 1239.2858
$ws.Range("J71").Value = 2430.2222
$ws.Range("K71").Value = 11153.5722
$ws.Range("L71").Value = 21871.9998
$ws.Range("M71").Value = -7097.572200000001
$ws.Range("N71").Value = -29983.9998
$ws.Range("H135").Value = 866.61536
$ws.Range("J135").Value = 1153.125
$ws.Range("L135").Value = 10378.125
$ws.Range("N135").Value = -15448.125

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 11868.714
$ws.Range("I132").Value = 12617.5
$ws.Range("K132").Value = 37852.5
$ws.Range("M132").Value = -35322.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1440.3125
$ws.Range("I16").Value = 788.72
$ws.Range("K16").Value = 788.72
$ws.Range("M16").Value = -618.72
$ws.Range("H122").Value = 4419.5264
$ws.Range("I122").Value = 4560.75
$ws.Range("K122").Value = 13682.25
$ws.Range("M122").Value = -11232.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3015.75
$ws.Range("J132").Value = 4750
$ws.Range("L132").Value = 14250
$ws.Range("N132").Value = -19310
